# Remove the obsolete AEB start/end deceleration-delta parameters
# (START_DECEL_DELTA / END_DECEL_DELTA) from the "params" sheet.
# These lived on rows 6-7; deleting the whole rows shifts everything
# below them up by two rows, matching the published diff.

$wb = $excel.ActiveWorkbook
$params = $wb.Worksheets.Item("params")
$params.Rows("6:7").Delete()

# Re-create the author's final UI state: selection on graphSpec was left
# at G29 when the user switched away from it, and "params" (now the
# active tab) ends up selected at D33.
$graphSpec = $wb.Worksheets.Item("graphSpec")
$graphSpec.Activate()
$graphSpec.Range("G29").Select()

$params.Activate()
$params.Range("D33").Select()
